# Update "Resumo de Inscricoes" figures on the "Inscricoes" sheet.
# These numbers reflect updated enrollment counts (Inscritos/Pagos/
# Inscricoes homologadas) for a number of process/campus rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4
$ws.Range("E4").Value = 40

# Row 11
$ws.Range("E11").Value = 291
$ws.Range("F11").Value = 161
$ws.Range("H11").Value = 161

# Row 12
$ws.Range("E12").Value = 417
$ws.Range("F12").Value = 224
$ws.Range("H12").Value = 224

# Row 13
$ws.Range("E13").Value = 110

# Row 15
$ws.Range("E15").Value = 138
$ws.Range("F15").Value = 55
$ws.Range("H15").Value = 55

# Row 17
$ws.Range("E17").Value = 81

# Row 22
$ws.Range("E22").Value = 148

# Row 24
$ws.Range("E24").Value = 182

# Row 25
$ws.Range("E25").Value = 224

# Row 27
$ws.Range("E27").Value = 286

# Row 28
$ws.Range("E28").Value = 170
$ws.Range("F28").Value = 59
$ws.Range("H28").Value = 59

# Row 32
$ws.Range("E32").Value = 164
$ws.Range("F32").Value = 98
$ws.Range("H32").Value = 98

# Row 33
$ws.Range("E33").Value = 252

# Row 34
$ws.Range("E34").Value = 185
$ws.Range("F34").Value = 114
$ws.Range("H34").Value = 114

# Row 35
$ws.Range("E35").Value = 123

# Row 36
$ws.Range("E36").Value = 58
$ws.Range("F36").Value = 36
$ws.Range("H36").Value = 36

# Row 41
$ws.Range("E41").Value = 338

# Row 42
$ws.Range("E42").Value = 312
$ws.Range("F42").Value = 166
$ws.Range("H42").Value = 166

# Row 45
$ws.Range("E45").Value = 123
$ws.Range("F45").Value = 61
$ws.Range("H45").Value = 61

# Row 46
$ws.Range("E46").Value = 272
$ws.Range("F46").Value = 148
$ws.Range("H46").Value = 148

# Row 47
$ws.Range("E47").Value = 385

# Row 49
$ws.Range("E49").Value = 257

# Row 50
$ws.Range("E50").Value = 221

# Row 51
$ws.Range("E51").Value = 209

$wb.Save()
